$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -13.543
$ws.Range("C4").Value = -12.73
$ws.Range("A11").Value = -21.822
$ws.Range("A12").Value = -21.721
$ws.Range("C14").Value = -12.083
$ws.Range("A15").Value = -22.006
$ws.Range("C26").Value = -12.805
$ws.Range("A27").Value = -21.80099999999999
$ws.Range("A28").Value = -21.779
$ws.Range("A31").Value = -22.024
$ws.Range("C31").Value = -13.27
$ws.Range("A32").Value = -21.697
$ws.Range("C35").Value = -13.176
$ws.Range("A36").Value = -20.275
$ws.Range("C37").Value = -13.293
$ws.Range("A38").Value = -19.69600000000001
$ws.Range("C39").Value = -12.967
$ws.Range("C40").Value = -12.51
$ws.Range("C45").Value = -12.477
$ws.Range("A46").Value = -21.838
$ws.Range("C52").Value = -11.587
$ws.Range("A54").Value = -22.15
$ws.Range("A55").Value = -22.151
$ws.Range("A56").Value = -21.988
$ws.Range("C57").Value = -13.49
$ws.Range("A67").Value = -21.519
$ws.Range("A69").Value = -21.721
$ws.Range("A72").Value = -21.481
$ws.Range("A73").Value = -19.994
$ws.Range("C81").Value = -13.192
$ws.Range("A83").Value = -21.658
$ws.Range("C83").Value = -12.789
$ws.Range("A86").Value = -22.257
$ws.Range("A91").Value = -21.522
$ws.Range("A93").Value = -21.421
$ws.Range("A99").Value = -19.861
$ws.Range("C100").Value = -13.018
$ws.Range("C102").Value = -13.221
